$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "EJBServer\components\MOLSA_ar\data\initial\WORKQUEUE.dmx`nEJBServer\components\MOLSA_ar\codetable\CT_ProgramWithdrawalRequestReason.ctx`nEJBServer\components\MOLSA\message\MOLSANotification.xml`nEJBServer\components\MOLSA\codetable\CT_ProgramWithdrawalRequestReason.ctx`nEJBServer\components\MOLSA\data\initial\clob\MOLSAIntakeConfiguration.xml`nEJBServer\components\MOLSA\data\initial\ALLOCATIONTARGET.dmx`nEJBServer\components\MOLSA\data\initial\ALLOCATIONTARGETITEM.dmx`nEJBServer\components\MOLSA\data\initial\MILESTONECONFIGURATION.dmx`nEJBServer\components\MOLSA\data\initial\WORKQUEUE.dmx`nEJBServer\components\MOLSA\events\handler_config.xml`nEJBServer\components\MOLSA\source\curam\molsa\constants\impl\MOLSAConstants.java`nEJBServer\components\MOLSA\source\curam\molsa\creoleprogramrecommendation\sl\event\impl\MOLSAApplicationDenialHandler.java`nEJBServer\components\MOLSA\source\curam\molsa\creoleprogramrecommendation\sl\event\impl\MolsaStatusChangedHandler.java`nEJBServer\components\MOLSA\workflow\MOLSAApplicationWithdrawTask_v1.xml"

$ws.Range("C9").Value = $newText
$ws.Range("C9").WrapText = $true

$ws.Range("E9").Value = "Smitha"
$ws.Range("F9").Value = "Yes"

$ws.Rows.Item(9).RowHeight = 240

$ws.Activate()
$ws.Range("D8").Select()
